$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.Value = "'" + $val
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '70.372.51'
Set-TextValue 'E2' '  +0.78%  '

Set-TextValue 'D3' '3.621.03'
Set-TextValue 'E3' '  +2.71%  '

Set-TextValue 'E4' '  +0.09%  '

Set-TextValue 'D5' '601.94'
Set-TextValue 'E5' '  -0.60%  '

Set-TextValue 'D6' '196.90'
Set-TextValue 'E6' '  +0.07%  '

Set-TextValue 'E7' '  -0.68%  '

Set-TextValue 'E8' '  +0.04%  '

Set-TextValue 'D9' '0.212'
Set-TextValue 'E9' '  +7.07%  '

Set-TextValue 'D10' '0.647'
Set-TextValue 'E10' '  -0.14%  '

Set-TextValue 'D11' '53.27'
Set-TextValue 'E11' '  -0.78%  '

Set-TextValue 'E12' '  +1.06%  '

Set-TextValue 'E13' '  +0.54%  '

Set-TextValue 'D14' '4.195.25'
Set-TextValue 'E14' '  +2.69%  '

Set-TextValue 'D15' '604.93'
Set-TextValue 'E15' '  +1.40%  '

Set-TextValue 'D16' '12.90'
Set-TextValue 'E16' '  +1.45%  '

Set-TextValue 'D17' '70.455.20'
Set-TextValue 'E17' '  +0.62%  '

Set-TextValue 'B18' 'WrappedEther'
Set-TextValue 'C18' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D18' '3.615.44'
Set-TextValue 'E18' '  +2.37%  '

Set-TextValue 'B19' 'Chainlink'
Set-TextValue 'C19' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D19' '19.04'
Set-TextValue 'E19' '  -0.31%  '

Set-TextValue 'E20' '  +1.41%  '

Set-TextValue 'D21' '0.999'
Set-TextValue 'E21' '  +0.80%  '

Set-TextValue 'D22' '18.15'

Set-TextValue 'D23' '5.19'
Set-TextValue 'E23' '  -1.45%  '

Set-TextValue 'D24' '103.03'
Set-TextValue 'E24' '  +1.44%  '

Set-TextValue 'D25' '4.61'
Set-TextValue 'E25' '  -0.11%  '

Set-TextValue 'D26' '2.99'
Set-TextValue 'E26' '  -6.70%  '

Set-TextValue 'D27' '10.64'
Set-TextValue 'E27' '  -2.41%  '

Set-TextValue 'D28' '9.70'
Set-TextValue 'E28' '  +1.22%  '

Set-TextValue 'E29' '  +1.38%  '

Set-TextValue 'D30' '4.68'
Set-TextValue 'E30' '  +7.48%  '

Set-TextValue 'E31' '  +2.86%  '

Set-TextValue 'E32' '  -1.23%  '

Set-TextValue 'E33' '  +1.32%  '

Set-TextValue 'D34' '63.35'
Set-TextValue 'E34' '  +0.41%  '

Set-TextValue 'D35' '0.0₃0890'
Set-TextValue 'E35' '  +3.69%  '

Set-TextValue 'D36' '3.924.53'
Set-TextValue 'E36' '  +5.96%  '

Set-TextValue 'E37' '  +0.19%  '

Set-TextValue 'B38' 'Bittensor'
Set-TextValue 'C38' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D38' '521.44'
Set-TextValue 'E38' '  +6.33%  '

Set-TextValue 'B39' 'Fetch.AI'
Set-TextValue 'C39' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D39' '3.06'
Set-TextValue 'E39' '  +0.19%  '

Set-TextValue 'D40' '36.80'
Set-TextValue 'E40' '  +0.50%  '

Set-TextValue 'E41' '  -0.85%  '

Set-TextValue 'E42' '  -2.34%  '

Set-TextValue 'D43' '0.137'
Set-TextValue 'E43' '  +3.18%  '

Set-TextValue 'D44' '0.0462'
Set-TextValue 'E44' '  +1.84%  '

Set-TextValue 'D45' '3.51'
Set-TextValue 'E45' '  +5.79%  '

Set-TextValue 'E46' '  +2.45%  '

Set-TextValue 'E47' '  -0.09%  '

Set-TextValue 'E48' '  +0.06%  '

Set-TextValue 'E49' '  -0.29%  '

Set-TextValue 'D50' '0.000250'
Set-TextValue 'E50' '  -0.46%  '

Set-TextValue 'E51' '  +0.69%  '
